$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Insert a new slide before the existing slide 2.
#
#    We duplicate slide 2 (which holds the picture) so the duplicate keeps
#    the picture and lands right after the original (i.e. becomes the new
#    slide 3), then we strip all shapes from the original slide 2 so it
#    becomes the new, empty, slide 2.
# ---------------------------------------------------------------------------
$p.Slides.Item(2).Duplicate() | Out-Null

$original = $p.Slides.Item(2)
for ($i = $original.Shapes.Count; $i -ge 1; $i--) {
    $original.Shapes.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# 2) Bump the "updated" date placeholder text from 2023/7/16 to 2023/7/17
#    on every slide layout and on the slide master.
# ---------------------------------------------------------------------------
$oldDate = "2023/7/16"
$newDate = "2023/7/17"

$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shp = $master.Shapes.Item($j)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
